$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Implementacion de intefaz punto" ---
$ws.Range("A6").Value = "Implementacion de intefaz punto"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 0.0048611111111111112
$ws.Range("E6").Value = 0.34166666666666662
$ws.Range("F6").Value = 0.34722222222222227
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# --- Row 7: "Implementacion de Punto2D" ---
$ws.Range("A7").Value = "Implementacion de Punto2D"
$ws.Range("B7").Value = 80
$ws.Range("C7").Value = 107
$ws.Range("D7").Value = 0.013888888888888888
$ws.Range("E7").Value = 0.34791666666666665
$ws.Range("F7").Value = 0.36041666666666666
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

# Row 7's label cell picks up the (border-less-top) formatting that was
# used for row 6 in the source workbook - match it by copying formats.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 8: "Implementacion de Punto3D" ---
$ws.Range("A8").Value = "Implementacion de Punto3D"
$ws.Range("B8").Value = 80
$ws.Range("D8").Value = 0.013888888888888888
$ws.Range("E8").Value = 0.36458333333333331
$ws.Range("F8").Value = 0.38194444444444442
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0

# The engine's incremental recalculation can leave a couple of
# two-hop-deep dependents (B12, B16) holding a value computed from an
# intermediate state seen earlier in this script. Re-asserting their
# formula (identical text) forces them to pick up the final values of
# their precedents without altering any other cell.
$b12 = $ws.Range("B12")
$b12.Formula = $b12.Formula
$b16 = $ws.Range("B16")
$b16.Formula = $b16.Formula

# Keep selection consistent with the edited workbook
$ws.Range("A9").Select()

$wb.Save()
